# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 325
$wsExpo.Range("F4").Value = 415
$wsExpo.Range("F5").Value = 1655
$wsExpo.Range("F7").Value = 2159
$wsExpo.Range("F8").Value = 4
$wsExpo.Range("F11").Value = 4822
$wsExpo.Range("F15").Value = 223
$wsExpo.Range("F16").Value = 27
$wsExpo.Range("F17").Value = 169
$wsExpo.Range("F18").Value = 34
$wsExpo.Range("F20").Value = 114
$wsExpo.Range("F21").Value = 3754
$wsExpo.Range("F22").Value = 691
$wsExpo.Range("F23").Value = 622
$wsExpo.Range("F31").Value = 570
$wsExpo.Range("F33").Value = 21
$wsExpo.Range("F34").Value = 872
$wsExpo.Range("F35").Value = 2376

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 325
$wsAll.Range("F4").Value = 415
$wsAll.Range("F5").Value = 1655
$wsAll.Range("F7").Value = 2159
$wsAll.Range("F8").Value = 4
$wsAll.Range("F11").Value = 4822
$wsAll.Range("F15").Value = 223
$wsAll.Range("F16").Value = 27
$wsAll.Range("F17").Value = 169
$wsAll.Range("F18").Value = 34
$wsAll.Range("F20").Value = 114
$wsAll.Range("F21").Value = 3754
$wsAll.Range("F22").Value = 691
$wsAll.Range("F23").Value = 622
$wsAll.Range("F31").Value = 570
$wsAll.Range("F34").Value = 21
$wsAll.Range("F35").Value = 872
$wsAll.Range("F36").Value = 2376
